# Normalize the "loc" column (B) from bracketed "[lat, lon]" text to a
# plain "lat,lon" string (drop the square brackets and the space after
# the comma) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value()
    if ($val -ne $null) {
        $trimmed = $val.Trim()
        if ($trimmed.StartsWith("[") -and $trimmed.EndsWith("]")) {
            $inner = $trimmed.Substring(1, $trimmed.Length - 2)
            $parts = $inner.Split(",")
            if ($parts.Length -eq 2) {
                $newVal = $parts[0].Trim() + "," + $parts[1].Trim()
                $cell.Value = $newVal
            }
        }
    }
}
